$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Agregado de etiquetas 306-HE PALETA BLOQ 4C - jabat
# Row 16 (4C / Paleta bloq) gets its label number filled in and the
# translation status flips from CERO (missing) to OK (labeled).
$ws.Range("M16").Value = 306
$ws.Range("N16").Value = "OK"

# The row grows taller to accommodate the extra content once filled in.
$ws.Rows.Item(16).RowHeight = 82.5

# Scroll the frozen view down a bit and leave the selection on the next
# row, matching where the user continued working after this entry.
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("N17").Select()

$wb.Save()
